# Update "Top Gainers" sheet data to reflect the 2025-10-30 13:03 refresh.
# The leaderboard is re-sorted by the "Latest" column, so some rows shift
# position (stock name changes) while others just get refreshed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top Gainers")

# Rows whose Stock name changed (rank shuffled in) -- update B:E (Stock, Latest, Weekly, Monthly)
$fullRows = @(
    @(10, "MCLOUD", 6.0394, 5.0231, -23.1658),
    @(11, "SHREEJISPG", 6.0362, 10.3075, 11.4016),
    @(12, "INDIACEM", 5.9537, 6.2388, 7.9818),
    @(13, "NETWEB", 5.8789, 11.778, 13.8797),
    @(14, "BLSE", 5.7958, 4.7583, -1.4574),
    @(15, "PDSL", 5.6805, 9.028700000000001, 15.1972),
    @(16, "VENKEYS", 5.3483, 5.9772, 3.6721),
    @(17, "POLICYBZR", 5.1526, 7.5043, 6.4747),
    @(23, "DBL", 4.8482, 5.9569, 7.0552),
    @(24, "BHEL", 4.8128, 11.3709, 7.8633),
    @(25, "SHRINGARMS", 4.7716, 5.9968, 26.4721),
    @(26, "RAMASTEEL", 4.6185, 4.5135, 6.11),
    @(28, "SHANTIGOLD", 4.3665, 11.7123, 4.2592),
    @(29, "MTARTECH", 4.2966, 8.449299999999999, 32.5537),
    @(32, "IVALUE", 4.1931, 7.4779, 0.1971),
    @(33, "HIRECT", 4.1509, 11.6032, 9.8894),
    @(35, "SKYGOLD", 3.9622, -0.6018, 38.0437),
    @(36, "FIVESTAR", 3.8427, 17.0019, 17.0891),
    @(43, "JKTYRE", 3.3597, 6.4207, 22.5779),
    @(44, "RELTD", 3.29, 10.0824, -1.4497),
    @(45, "BLISSGVS", 3.2555, 2.6077, 2.9778),
    @(46, "PSPPROJECT", 3.2262, 17.3483, 23.8),
    @(47, "BGRENERGY", 3.2153, -6.0917, 74.82810000000001),
    @(48, "CANBK", 3.2073, 5.7279, 7.4287),
    @(49, "MRPL", 3.1991, 13.22, 23.8949),
    @(50, "BLS", 3.1893, 0.1414, -1.117),
    @(51, "GRMOVER", 3.1273, 3.3019, 19.0051),
    @(52, "PFOCUS", 3.1013, 0.3922, 1.8472),
    @(53, "NEULANDLAB", 3.0897, -1.2805, 8.743600000000001),
    @(54, "GMMPFAUDLR", 3.0423, 7.508, 19.8547),
    @(55, "RSYSTEMS", 3.0405, 4.1972, 6.5369),
    @(56, "SUNDROP", 3.0242, 2.9169, 1.0161),
    @(57, "GANESHCP", 3.0099, 2.4837, 2.0344),
    @(58, "ASHOKA", 3.0022, 4.5251, 7.156),
    @(59, "OIL", 2.7467, 2.9917, 4.4103),
    @(60, "VOLTAMP", 2.7398, 2.6188, 2.3564),
    @(61, "BPCL", 2.7004, 8.1858, 5.2554),
    @(62, "POWERINDIA", 2.6999, 7.0794, -0.2611),
    @(64, "IIFL", 2.6099, 9.565099999999999, 18.7542),
    @(66, "MFSL", 2.5876, 2.6416, -1.1295),
    @(67, "INDORAMA", 2.5364, 5.0007, 16.3799),
    @(68, "OBEROIRLTY", 2.5096, 3.3237, 11.017),
    @(69, "SDBL", 2.3981, 0.9456, 6.5266),
    @(70, "BLUEDART", 2.3661, 20.8786, 17.9642),
    @(71, "CIFL", 2.3553, 1.9108, 1.8519),
    @(72, "CARYSIL", 2.3377, 1.8236, 10.6831),
    @(73, "FEDFINA", 2.323, 3.4229, -5.2972),
    @(74, "JKLAKSHMI", 2.3142, 4.298, 1.3238),
    @(75, "TDPOWERSYS", 2.3089, 6.6355, 15.867),
    @(76, "SPANDANA", 2.2922, 3.7816, 2.4024)
)

# Rows whose Stock name stayed the same -- update only C:E (Latest, Weekly, Monthly)
$partialRows = @(
    @(2, 11.1784, 19.9849, 27.1428),
    @(3, 10.4629, 10.6693, 24.73),
    @(5, 9.7902, 6.8546, -7.9179),
    @(6, 8.0991, 15.7197, 28.903),
    @(7, 8.028600000000001, 10.3123, 26.5542),
    @(8, 7.4225, 12.6173, 14.7796),
    @(9, 6.5158, 3.0391, 15.6218),
    @(18, 5.1436, 3.754, 9.067399999999999),
    @(22, 4.8547, 10.7577, 27.8008),
    @(30, 4.2915, 4.8396, -1.4512),
    @(31, 4.2513, 10.4315, 15.9169),
    @(37, 3.8123, 12.1671, 11.7024),
    @(38, 3.6387, 3.5381, -3.3918),
    @(39, 3.6198, 5.4469, 2.791),
    @(42, 3.3939, 2.003, 1.3365),
    @(63, 2.6508, 10.654, 34.0162)
)

foreach ($row in $fullRows) {
    $r = [string]$row[0]
    $arr = New-Object 'object[,]' 1,4
    $arr[0,0] = $row[1]
    $arr[0,1] = $row[2]
    $arr[0,2] = $row[3]
    $arr[0,3] = $row[4]
    $addr = "B" + $r + ":E" + $r
    $ws.Range($addr).Value = $arr
}

foreach ($row in $partialRows) {
    $r = [string]$row[0]
    $arr = New-Object 'object[,]' 1,3
    $arr[0,0] = $row[1]
    $arr[0,1] = $row[2]
    $arr[0,2] = $row[3]
    $addr = "C" + $r + ":E" + $r
    $ws.Range($addr).Value = $arr
}
